# Update the "survey" sheet question labels to include numbered prefixes,
# and change the black-fly month question from a single-select to a
# multi-select. Also bump the form title/id from V2 to V3 on the
# "settings" sheet.

$wb = $excel.ActiveWorkbook

$survey = $wb.Worksheets.Item("survey")

# Prefix each question label (column C = label::English) with its
# question number, preserving the original (non-breaking-space) whitespace.
$survey.Range("C2").Value = "1. Enter recorder ID   "
$survey.Range("C3").Value = "2. Select district   "
$survey.Range("C5").Value = "3. Select community name"
$survey.Range("C7").Value = "4. Does this person consent to take part in the study?   "
$survey.Range("C8").Value = "5. What gender is this person?   "
$survey.Range("C9").Value = "6. What age is this person?   "
$survey.Range("C10").Value = "7. What is their occupation?   "
$survey.Range("C11").Value = "8. How long has this person lived in this community?   "
$survey.Range("C12").Value = "9. Are black fly bites a problem in the community? "
$survey.Range("C13").Value = "10. Bites number    "
$survey.Range("C14").Value = "11. Black fly month   "
$survey.Range("C15").Value = "12. Black fly time of day  "
$survey.Range("C16").Value = "13. Onchocerciasis – do you know what onchocerciasis/river blindness is? "
$survey.Range("C17").Value = "14. Can you explain what onchocerciasis/river blindness is? "
$survey.Range("C18").Value = "15. Enter GPS  "

# Row 14 ("Black fly month") becomes a multi-select question.
$survey.Range("A14").Value = "select_multiple month"

# Bump the form version from V2 to V3.
$settings = $wb.Worksheets.Item("settings")
$settings.Range("A2").Value = "(2022 October) - 1. vector community questionnaire V3"
$settings.Range("B2").Value = "gn_oncho_bsa_1_vector_community_question_202210_v3"

# Leave the cursor where the author last left it before saving.
$survey.Activate()
$survey.Range("A13").Select()
